$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Mon, Jan 13 / Taxonomy of Questions): expand "Do Before Class" text
#     to add a 3rd bullet linking to the new taxonomy_of_questions notebook, and
#     grow the row height to fit the extra line.
$c3 = $ws.Cells.Item(3, 3)
$c3.Value = "'- Read and sign syllabus`n- Submit substantive interest survey`n- ``Taxonomy of Questions <taxonomy_of_questions.ipynb>```_"
$ws.Rows.Item(3).RowHeight = 57

# --- Row 6 (Wed, Jan 22 / Descriptive Questions): replace the old bullet with a
#     link to the new "Mere Description" reading, and grow the row height.
$c6 = $ws.Cells.Item(6, 3)
$c6.Value = "'- ``Mere Description <https://doi.org/10.1017/S0007123412000130>```_"
$ws.Rows.Item(6).RowHeight = 43

# --- Row 12 (Wed, Feb 12 / CI: Regression): prefix the Angrist & Pischke
#     citation with a leading bullet dash, preserving the bold run on "Pages 249- 268."
$c12 = $ws.Cells.Item(12, 3)
$c12.Value = "'- Angrist and Piscke, Pages 249- 268. "
$run1Len = ("- Angrist and Piscke, ").Length
$run2Len = ("Pages 249- 268. ").Length
$boldRun = $c12.Characters($run1Len + 1, $run2Len)
$boldRun.Font.Size = 12
$boldRun.Font.Name = "SFBX1200"

# --- Row 13 (Mon, Feb 17 / CI: Fixed Effects): add a 3rd bullet (Callahan reading)
#     to the "Do Before Class" notes, and grow the row height to fit it.
$c13 = $ws.Cells.Item(13, 3)
$c13.Value = "'- ``Fixed Effects v. Hierarchical Models <fixed_effects_v_hierarchical.ipynb>```_`n- ``Interpreting Indicator Vars <interpreting_indicator_vars.ipynb>```_`n- Callahan, pp. 72-89 `n"
$ws.Rows.Item(13).RowHeight = 153

# --- Sheet view: scroll back to the top and select C4 instead of the old C30 selection.
$ws.Range("C4").Select() | Out-Null
